$d = $word.ActiveDocument

# 1) "End game when hit "q"" -> "Regain hit points at start of turn"
$d.Content.Find.Execute("End game when hit “q”", $false, $false, $false, $false, $false, $true, 1, $false, "Regain hit points at start of turn", 2) | Out-Null

# 2) The paragraph that used to read "Reg" + bookmark + "ain hit points at start
#    of turn" (net text "Regain hit points at start of turn") -> "Scrolls".
#    Scope the Find to paragraph 3 specifically so we don't touch paragraph 2,
#    which now has the same resulting text after step 1.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Find.Execute("Regain hit points at start of turn", $false, $false, $false, $false, $false, $true, 1, $false, "Scrolls", 2) | Out-Null

# 3) "Cheat mode" -> "Sleep"
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Find.Execute("Cheat mode", $false, $false, $false, $false, $false, $true, 1, $false, "Sleep", 2) | Out-Null

# The original "_GoBack" bookmark (which used to sit between "Reg" and "ain...")
# was consumed by the step-2 replace. Recreate it as a zero-length bookmark
# right after "Sleep", matching the post-edit document.
#
# Note: adding a bookmark on a collapsed Range positioned exactly at a
# paragraph's own end (i.e. right before its paragraph mark) lands in the
# wrong place, so as a workaround we briefly insert a placeholder character
# after "Sleep", anchor the bookmark at the boundary (now a safe, non-edge
# position), and then remove the placeholder again.
$p4 = $d.Paragraphs.Item(4)
$r = $p4.Range
$r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark -> exactly "Sleep"
$splitPos = $r.End
$r.InsertAfter("X") | Out-Null

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$placeholder = $d.Range($splitPos, $splitPos + 1)
$placeholder.Delete() | Out-Null
